$wb = $excel.ActiveWorkbook

# --- Add the new "Setting" sheet, placed before the existing "Login" sheet ---
$loginSheet = $wb.Worksheets.Item("Login")
$settingSheet = $wb.Worksheets.Add($loginSheet)
$settingSheet.Name = "Setting"

# --- Populate the new sheet's data ---
$settingSheet.Range("A1").Value = "checkStoreDB"
$settingSheet.Range("B1").Value = "checkPaging"
$settingSheet.Range("A2").Value = "Yes"
$settingSheet.Range("B2").Value = "No"

# --- Style the new cells like the rest of the workbook (thin box border, style index 1) ---
$settingSheet.Range("A1:B2").Borders.LineStyle = 1
$settingSheet.Range("A1:B2").Borders.Weight = 2

# --- Column widths approximating the original author's "best fit" auto-sizing ---
$settingSheet.Columns.Item(1).ColumnWidth = 12.42
$settingSheet.Columns.Item(2).ColumnWidth = 10.92

# --- Data validation dropdown (Yes/No) on A2:B2 ---
$settingSheet.Range("A2:B2").Validation.Add(3, 1, 1, '"Yes, No"')

# --- Selection / active cell on the new sheet ---
$settingSheet.Range("B3").Select()

# --- Make the new sheet the active/selected tab ---
$settingSheet.Activate()
